# Updated cryptos list on Tue Oct 31 07:19:32 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of "Price" cells get new values that Excel would otherwise
# auto-detect as numbers (single decimal point, e.g. "35.99"). The source
# data keeps these as plain text (inline strings), so force a Text number
# format on just those cells before writing the new value, preserving the
# original text representation instead of letting COM coerce it to a double.
$textCells = @("D5","D6","D8","D9","D13","D18","D19","D21","D25","D26","D27","D32","D33","D36","D41","D42","D43","D45","D46","D47","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.355.73"

$ws.Range("D3").Value = "1.801.39"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("D5").Value = "227.19"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").Value = "0.575"
$ws.Range("E6").Value = "  +3.92%  "

$ws.Range("D8").Value = "35.99"
$ws.Range("E8").Value = "  +9.77%  "

$ws.Range("D9").Value = "0.299"
$ws.Range("E9").Value = "  +1.29%  "

$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("E11").Value = "  +1.99%  "

$ws.Range("D12").Value = "2.061.36"
$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("D13").Value = "11.66"
$ws.Range("E13").Value = "  +5.37%  "

$ws.Range("D14").Value = "1.818.59"
$ws.Range("E14").Value = "  +1.42%  "

$ws.Range("E15").Value = "  +1.24%  "

$ws.Range("E16").Value = "  +4.97%  "

$ws.Range("D17").Value = "34.352.18"
$ws.Range("E17").Value = "  -0.07%  "

$ws.Range("D18").Value = "69.03"
$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("D19").Value = "245.05"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("D21").Value = "11.62"
$ws.Range("E21").Value = "  +3.13%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("E23").Value = "  -0.09%  "

$ws.Range("E24").Value = "  +3.87%  "

$ws.Range("D25").Value = "171.59"
$ws.Range("E25").Value = "  +3.00%  "

$ws.Range("D26").Value = "7.96"
$ws.Range("E26").Value = "  +8.99%  "

$ws.Range("D27").Value = "16.84"
$ws.Range("E27").Value = "  +2.01%  "

$ws.Range("E28").Value = "  +1.83%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("E31").Value = "  +0.91%  "

# Rows 32/33 swap: PancakeSwap <-> Filecoin change places in the ranking.
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "3.83"
$ws.Range("E32").Value = "  +0.59%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "1.24"
$ws.Range("E33").Value = "  +1.05%  "

$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("D35").Value = "1.398.29"
$ws.Range("E35").Value = "  -0.28%  "

$ws.Range("D36").Value = "0.669"
$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("E39").Value = "  -0.28%  "

$ws.Range("E40").Value = "  +10.38%  "

$ws.Range("D41").Value = "0.962"
$ws.Range("E41").Value = "  +2.60%  "

$ws.Range("D42").Value = "82.70"
$ws.Range("E42").Value = "  -3.26%  "

$ws.Range("D43").Value = "2.82"
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("D45").Value = "13.44"
$ws.Range("E45").Value = "  -2.64%  "

$ws.Range("D46").Value = "0.0507"
$ws.Range("E46").Value = "  -3.63%  "

$ws.Range("D47").Value = "6.04"
$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("D48").Value = "1.961.82"
$ws.Range("E48").Value = "  +0.62%  "

$ws.Range("D49").Value = "104.23"
$ws.Range("E49").Value = "  -0.76%  "

$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("E51").Value = "  +0.07%  "
